$d = $word.ActiveDocument

$replacements = @(
    @("2025-08-16 Saturday", "2025-08-17 Sunday"),
    @("98×14=1372", "97×96=9312"),
    @("50×49=2450", "62×17=1054"),
    @("71×69=4899", "98×55=5390"),
    @("45×99=4455", "59×93=5487"),
    @("46×84=3864", "17×81=1377"),
    @("78×63=4914", "19×36=684"),
    @("45×95=4275", "86×27=2322"),
    @("80×15=1200", "17×69=1173"),
    @("76×30=2280", "60×39=2340"),
    @("74×36=2664", "67×29=1943"),
    @("19×74=1406", "28×49=1372"),
    @("49×53=2597", "77×81=6237"),
    @("34×38=1292", "93×15=1395"),
    @("14×41=574", "57×48=2736"),
    @("52×78=4056", "80×86=6880"),
    @("72×95=6840", "77×61=4697"),
    @("96×80=7680", "22×19=418"),
    @("55×99=5445", "37×70=2590"),
    @("28×72=2016", "43×88=3784"),
    @("60×22=1320", "45×63=2835"),
    @("47×25=1175", "48×40=1920"),
    @("99×95=9405", "51×62=3162"),
    @("69×81=5589", "84×88=7392"),
    @("40×39=1560", "65×85=5525"),
    @("16×26=416", "89×51=4539")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
